# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values (GitHub Actions cron refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text entry so numeric-looking strings
# (e.g. '582.75') are NOT reinterpreted as numbers - matches the
# original sheet, where every Price/Volume cell is stored as text.

$ws.Range("D2").Value = "'68.040.39"
$ws.Range("E2").Value = "'  +0.67%  "
$ws.Range("D3").Value = "'3.262.88"
$ws.Range("E3").Value = "'  +0.80%  "
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("D5").Value = "'582.75"
$ws.Range("E5").Value = "'  +0.71%  "
$ws.Range("D6").Value = "'185.32"
$ws.Range("E6").Value = "'  +2.08%  "
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E8").Value = "'  +1.00%  "
$ws.Range("E9").Value = "'  -2.87%  "
$ws.Range("D10").Value = "'6.61"
$ws.Range("E10").Value = "'  +0.31%  "
$ws.Range("E11").Value = "'  -2.29%  "
$ws.Range("D12").Value = "'3.831.18"
$ws.Range("E12").Value = "'  +1.05%  "
$ws.Range("E13").Value = "'  +1.60%  "
$ws.Range("E14").Value = "'  -1.94%  "
$ws.Range("D15").Value = "'68.068.18"
$ws.Range("E15").Value = "'  +0.70%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("E16").Value = "'  -1.73%  "
$ws.Range("D17").Value = "'3.287.58"
$ws.Range("E17").Value = "'  +3.23%  "
$ws.Range("D18").Value = "'5.70"
$ws.Range("E18").Value = "'  -1.55%  "
$ws.Range("E19").Value = "'  -1.12%  "
$ws.Range("D20").Value = "'415.79"
$ws.Range("E20").Value = "'  +6.77%  "
$ws.Range("E21").Value = "'  -1.63%  "
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("D23").Value = "'71.36"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E24").Value = "'  -1.48%  "
$ws.Range("E25").Value = "'  -1.39%  "
$ws.Range("E26").Value = "'  -0.56%  "
$ws.Range("E27").Value = "'  -1.21%  "
$ws.Range("E28").Value = "'  +0.78%  "
$ws.Range("E29").Value = "'  -1.11%  "
$ws.Range("D30").Value = "'22.63"
$ws.Range("E30").Value = "'  -1.41%  "
$ws.Range("E31").Value = "'  -3.04%  "
$ws.Range("E32").Value = "'  +0.01%  "
$ws.Range("D33").Value = "'6.86"
$ws.Range("E33").Value = "'  -3.14%  "
$ws.Range("E34").Value = "'  -2.54%  "
$ws.Range("D35").Value = "'162.66"
$ws.Range("E35").Value = "'  -0.60%  "
$ws.Range("E36").Value = "'  -2.53%  "
$ws.Range("E37").Value = "'  -0.43%  "
$ws.Range("D38").Value = "'27.01"
$ws.Range("E38").Value = "'  +1.69%  "
$ws.Range("D39").Value = "'0.796"
$ws.Range("E39").Value = "'  -2.22%  "
$ws.Range("E40").Value = "'  -2.91%  "
$ws.Range("D41").Value = "'6.35"
$ws.Range("E41").Value = "'  -2.08%  "
$ws.Range("D42").Value = "'2.637.24"
$ws.Range("E42").Value = "'  +0.84%  "
$ws.Range("D43").Value = "'40.78"
$ws.Range("E43").Value = "'  -1.39%  "
$ws.Range("E44").Value = "'  -2.66%  "
$ws.Range("E45").Value = "'  -0.92%  "
$ws.Range("D46").Value = "'337.17"
$ws.Range("E46").Value = "'  -0.95%  "
$ws.Range("D47").Value = "'24.35"
$ws.Range("E47").Value = "'  -1.38%  "
$ws.Range("D48").Value = "'0.0274"
$ws.Range("E48").Value = "'  -2.01%  "
$ws.Range("E49").Value = "'  -0.54%  "
$ws.Range("D50").Value = "'0.977"
$ws.Range("E50").Value = "'  +0.47%  "
$ws.Range("E51").Value = "'  -1.35%  "
